# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (avoids Excel auto-converting numeric-looking strings to numbers,
# which would strip trailing zeros / use scientific notation).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.388.37'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.687.43'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '679.33'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = '159.30'
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('D10').Value = '7.09'
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').Value = '4.311.09'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '32.41'
$ws.Range('E14').Value = '  -3.25%  '
$ws.Range('D15').Value = '3.681.50'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').Value = '69.303.33'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').Value = '16.05'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '6.43'
$ws.Range('E19').Value = '  -2.86%  '
$ws.Range('D20').Value = '468.50'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').Value = '10.01'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  -2.26%  '
$ws.Range('D23').Value = '79.92'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '3.832.90'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D26').Value = '0.0000123'
$ws.Range('E26').Value = '  -5.57%  '
$ws.Range('D27').Value = '10.93'
$ws.Range('E27').Value = '  -4.45%  '
$ws.Range('D28').Value = '9.12'
$ws.Range('E28').Value = '  -4.40%  '
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('E30').Value = '  -3.24%  '
$ws.Range('E31').Value = '  -3.28%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '26.96'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = '3.676.91'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -4.87%  '
$ws.Range('D37').Value = '8.27'
$ws.Range('E37').Value = '  -2.72%  '
$ws.Range('D38').Value = '6.23'
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '2.25'
$ws.Range('E40').Value = '  -3.05%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('D43').Value = '170.41'
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('D44').Value = '0.943'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '28.11'
$ws.Range('E47').Value = '  -6.09%  '
$ws.Range('B48').Value = 'FLOKI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D48').Value = '0.000278'
$ws.Range('E48').Value = '  -3.00%  '
$ws.Range('E49').Value = '  -3.34%  '
$ws.Range('D50').Value = '1.29'
$ws.Range('E50').Value = '  -5.32%  '
$ws.Range('D51').Value = '7.80'
$ws.Range('E51').Value = '  -2.63%  '
